# ComplianceReports workbook update:
# Mark previously "Not verified / Need Verification" test cases as fixed
# (Status P -> A, Comment prefixed with "Fixed: ") now that SSRS PDF and
# Views verifications were added to ComplianceReportsPageTest2, and move
# the active tab/selection to reflect the sheets that were worked on.

$wb = $excel.ActiveWorkbook

function Mark-Fixed($ws, $row) {
    $ws.Cells.Item($row, 2).Value = "A"
    $commentCell = $ws.Cells.Item($row, 3)
    $orig = $commentCell.Value2
    $commentCell.Value = "Fixed: " + $orig
}

# --- Sheet "CRPT" ---
$crpt = $wb.Worksheets.Item("CRPT")
$crpt.Activate()
Mark-Fixed $crpt 9
$crpt.Range("C17").Select()

# --- Sheet "CRPTEthane" ---
$ethane = $wb.Worksheets.Item("CRPTEthane")
$ethane.Activate()
10..21 | ForEach-Object { Mark-Fixed $ethane $_ }
$ethane.Range("C10:C21").Select()

# --- Sheet "CRPT-2" ---
$crpt2 = $wb.Worksheets.Item("CRPT-2")
$crpt2.Activate()
@(3,5,6,7,8,9,10,25) | ForEach-Object { Mark-Fixed $crpt2 $_ }
$crpt2.Range("A10").Select()
